# arreglo menu lateral - datos financieros
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update claim numbers (NroSiniestro column) in rows 2 and 3.
# A leading apostrophe forces these numeric-looking claim ids to be
# stored as text, matching their original type (no format change).
# Written E3-then-E2 so the shared-string table keeps the same slot
# ordering as before the edit.
$ws.Range("E3").Value = "'1120194100412"
$ws.Range("E2").Value = "'1220194200667"

# Update the active selection stored with the sheet view.
$ws.Range("I6:I7").Select()
